$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (Changed) date column C for rows 2-31 from 45183 to 45184
$ws.Range("C2:C31").Value = 45184
